$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "0.530", "322.80", "1.00") are preserved verbatim rather than
# being normalized as numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '61.983.29'
$ws.Range("E2").Value = '  +1.28%  '

$ws.Range("D3").Value = '2.411.29'
$ws.Range("E3").Value = '  +1.39%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Value = '556.78'
$ws.Range("E5").Value = '  +1.34%  '

$ws.Range("D6").Value = '142.46'
$ws.Range("E6").Value = '  +3.06%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").Value = '0.530'
$ws.Range("E8").Value = '  +0.47%  '

$ws.Range("D9").Value = '2.404.63'
$ws.Range("E9").Value = '  +1.08%  '

$ws.Range("E10").Value = '  +0.88%  '

$ws.Range("E11").Value = '  -1.13%  '

$ws.Range("E12").Value = '  +0.36%  '

$ws.Range("D13").Value = '0.352'

$ws.Range("D14").Value = '25.98'
$ws.Range("E14").Value = '  +3.53%  '

$ws.Range("E15").Value = '  +3.94%  '

$ws.Range("D16").Value = '2.839.54'
$ws.Range("E16").Value = '  +1.75%  '

$ws.Range("D17").Value = '61.777.44'
$ws.Range("E17").Value = '  +1.08%  '

$ws.Range("D18").Value = '2.406.22'
$ws.Range("E18").Value = '  +0.71%  '

$ws.Range("D19").Value = '11.12'
$ws.Range("E19").Value = '  +2.78%  '

$ws.Range("D20").Value = '4.18'
$ws.Range("E20").Value = '  +0.50%  '

$ws.Range("D21").Value = '322.80'
$ws.Range("E21").Value = '  +0.43%  '

$ws.Range("D22").Value = '6.72'
$ws.Range("E22").Value = '  +0.11%  '

$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("D24").Value = '65.11'
$ws.Range("E24").Value = '  +1.32%  '

$ws.Range("D25").Value = '1.72'
$ws.Range("E25").Value = '  +1.84%  '

$ws.Range("E26").Value = '  +7.73%  '

$ws.Range("D27").Value = '580.82'
$ws.Range("E27").Value = '  +13.98%  '

$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.11%  '

$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '2.528.33'
$ws.Range("E29").Value = '  +1.64%  '

$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0929'
$ws.Range("E30").Value = '  +4.68%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '8.26'
$ws.Range("E31").Value = '  +1.17%  '

$ws.Range("D32").Value = '1.45'
$ws.Range("E32").Value = '  +3.95%  '

$ws.Range("E33").Value = '  -1.64%  '

$ws.Range("E34").Value = '  +2.29%  '

$ws.Range("D35").Value = '1.56'
$ws.Range("E35").Value = '  +2.36%  '

$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  +0.00%  '

$ws.Range("D37").Value = '5.66'
$ws.Range("E37").Value = '  +5.46%  '

$ws.Range("D38").Value = '4.75'
$ws.Range("E38").Value = '  +0.58%  '

$ws.Range("D39").Value = '0.383'
$ws.Range("E39").Value = '  +0.99%  '

$ws.Range("D40").Value = '150.97'
$ws.Range("E40").Value = '  +3.41%  '

$ws.Range("D41").Value = '18.66'
$ws.Range("E41").Value = '  +0.22%  '

$ws.Range("D42").Value = '1.83'
$ws.Range("E42").Value = '  -2.89%  '

$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.05%  '

$ws.Range("D44").Value = '2.30'
$ws.Range("E44").Value = '  +10.70%  '

$ws.Range("D45").Value = '150.71'
$ws.Range("E45").Value = '  +1.40%  '

$ws.Range("D46").Value = '3.64'
$ws.Range("E46").Value = '  +0.76%  '

$ws.Range("D47").Value = '0.0539'
$ws.Range("E47").Value = '  +3.22%  '

$ws.Range("D48").Value = '20.18'
$ws.Range("E48").Value = '  +3.97%  '

$ws.Range("D49").Value = '0.588'
$ws.Range("E49").Value = '  +1.96%  '

$ws.Range("D50").Value = '0.0923'
$ws.Range("E50").Value = '  +1.20%  '

$ws.Range("D51").Value = '0.0228'
$ws.Range("E51").Value = '  +1.77%  '
